$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell "time_taken" in F1, matching the style of the other
# header cells (bold font, border, centered alignment) by copying E1's format.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Timestamp values for each data row (2..24), column F, default style
# (matching the unstyled data cells in columns B-E).
$timestamps = @(
    "2021-10-05 13:39:41.737483",
    "2021-10-05 13:39:41.737496",
    "2021-10-05 13:39:41.737500",
    "2021-10-05 13:39:41.737503",
    "2021-10-05 13:39:41.737506",
    "2021-10-05 13:39:41.737510",
    "2021-10-05 13:39:41.737513",
    "2021-10-05 13:39:41.737516",
    "2021-10-05 13:39:41.737519",
    "2021-10-05 13:39:41.737522",
    "2021-10-05 13:39:41.737525",
    "2021-10-05 13:39:41.737528",
    "2021-10-05 13:39:41.737532",
    "2021-10-05 13:39:41.737534",
    "2021-10-05 13:39:41.737537",
    "2021-10-05 13:39:41.737540",
    "2021-10-05 13:39:41.737544",
    "2021-10-05 13:39:41.737547",
    "2021-10-05 13:39:41.737550",
    "2021-10-05 13:39:41.737553",
    "2021-10-05 13:39:41.737557",
    "2021-10-05 13:39:41.737560",
    "2021-10-05 13:39:41.737563"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
